# "Signed Off time sheets"
# Fill in the supervisor's initials on the "Week of" header row, and add the
# supervisor's signed-off name + date on the sign-off block at the bottom of
# the sheet (mirrors the existing employee sign-off row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supervisor sign-off block (row 27): name + sign-off date, mirroring the
# employee sign-off row above it (row 25).
$ws.Range("A27").Value = "Ankita Gangotra"
$ws.Range("D27").NumberFormat = "mm-dd-yy"
$ws.Range("D27").Value = (Get-Date -Year 2014 -Month 6 -Day 10 -Hour 0 -Minute 0 -Second 0)

# Supervisor initials next to "Week of:" (row 6, merged G6:I6)
$ws.Range("G6").Value = "A.G"

# Leave the selection on the cell that was last edited, as the saved workbook shows.
$ws.Range("D27:E27").Select()
